# Auto-generated script applying scheduled market-data refresh to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9983.666999999999
$ws.Range("I40").Value = 5451
$ws.Range("K40").Value = 5451
$ws.Range("M40").Value = -5276
$ws.Range("H86").Value = 3291748.5
$ws.Range("I86").Value = 2033.75
$ws.Range("K86").Value = 2033.75
$ws.Range("M86").Value = -910.75
$ws.Range("H89").Value = 3291748.5
$ws.Range("I89").Value = 2033.75
$ws.Range("K89").Value = 10168.75
$ws.Range("M89").Value = -4552.75
$ws.Range("H92").Value = 1019.73334
$ws.Range("I92").Value = 469.85715
$ws.Range("J92").Value = 1500.875
$ws.Range("K92").Value = 469.85715
$ws.Range("L92").Value = 1500.875
$ws.Range("M92").Value = 778.14285
$ws.Range("N92").Value = -3996.875
$ws.Range("H99").Value = 1471
$ws.Range("I99").Value = 177.5
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 532.5
$ws.Range("L99").Value = 6999.999899999999
$ws.Range("M99").Value = 965.5
$ws.Range("N99").Value = -9995.999899999999
$ws.Range("H100").Value = 6353.3335
$ws.Range("J100").Value = 17000.2
$ws.Range("L100").Value = 17000.2
$ws.Range("N100").Value = -18082.2
$ws.Range("H103").Value = 475.8889
$ws.Range("J103").Value = 497.57144
$ws.Range("L103").Value = 1492.71432
$ws.Range("N103").Value = -2664.71432
$ws.Range("H107").Value = 101396.3
$ws.Range("I107").Value = 144524.72
$ws.Range("K107").Value = 144524.72
$ws.Range("M107").Value = -142604.72
$ws.Range("H132").Value = 1752.6666
$ws.Range("I132").Value = 1157
$ws.Range("J132").Value = 4220.4287
$ws.Range("K132").Value = 3471
$ws.Range("L132").Value = 12661.2861
$ws.Range("M132").Value = -941
$ws.Range("N132").Value = -17721.2861
$ws.Range("H133").Value = 58747.645
$ws.Range("J133").Value = 58747.645
$ws.Range("L133").Value = 58747.645
$ws.Range("N133").Value = -68867.64499999999
$ws.Range("H135").Value = 986.1579
$ws.Range("I135").Value = 835.86664
$ws.Range("J135").Value = 1549.75
$ws.Range("K135").Value = 7522.79976
$ws.Range("L135").Value = 13947.75
$ws.Range("M135").Value = -4987.79976
$ws.Range("N135").Value = -19017.75
$ws.Range("H137").Value = 3503.3655
$ws.Range("I137").Value = 2169.4119
$ws.Range("J137").Value = 6023.0557
$ws.Range("K137").Value = 6508.2357
$ws.Range("L137").Value = 18069.1671
$ws.Range("M137").Value = -3958.2357
$ws.Range("N137").Value = -23169.1671
$ws.Range("H138").Value = 4170.873
$ws.Range("I138").Value = 2559.2856
$ws.Range("J138").Value = 5460.143
$ws.Range("K138").Value = 7677.8568
$ws.Range("L138").Value = 16380.429
$ws.Range("M138").Value = -2537.8568
$ws.Range("N138").Value = -26660.429
$ws.Range("H141").Value = 2729.3572
$ws.Range("I141").Value = 1428.2727
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 4284.8181
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = 895.1818999999996
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3954.1272
$ws.Range("I32").Value = 3429.78
$ws.Range("K32").Value = 3429.78
$ws.Range("M32").Value = -3142.78
$ws.Range("H45").Value = 4507.7393
$ws.Range("J45").Value = 5176.9287
$ws.Range("L45").Value = 5176.9287
$ws.Range("N45").Value = -5930.9287
$ws.Range("H61").Value = 2026.5807
$ws.Range("I61").Value = 1363.9259
$ws.Range("K61").Value = 1363.9259
$ws.Range("M61").Value = -1151.9259
$ws.Range("H74").Value = 1472.5
$ws.Range("I74").Value = 1141.32
$ws.Range("J74").Value = 4232.3335
$ws.Range("K74").Value = 1141.32
$ws.Range("L74").Value = 4232.3335
$ws.Range("M74").Value = -267.3199999999999
$ws.Range("N74").Value = -5980.3335
$ws.Range("H77").Value = 1472.5
$ws.Range("I77").Value = 1141.32
$ws.Range("J77").Value = 4232.3335
$ws.Range("K77").Value = 5706.599999999999
$ws.Range("L77").Value = 21161.6675
$ws.Range("M77").Value = -1338.599999999999
$ws.Range("N77").Value = -29897.6675
$ws.Range("H102").Value = 2105.8696
$ws.Range("I102").Value = 2030.2858
$ws.Range("K102").Value = 2030.2858
$ws.Range("M102").Value = -408.2858000000001
$ws.Range("H132").Value = 5258.909
$ws.Range("I132").Value = 3091.244
$ws.Range("K132").Value = 9273.732
$ws.Range("M132").Value = -6743.732
$ws.Range("H136").Value = 2026.5807
$ws.Range("I136").Value = 1363.9259
$ws.Range("K136").Value = 4091.7777
$ws.Range("M136").Value = -1541.7777
$ws.Range("H138").Value = 50750
$ws.Range("J138").Value = 50750
$ws.Range("L138").Value = 50750
$ws.Range("N138").Value = -61030
$ws.Range("H139").Value = 47849.8
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47849.8
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47849.8
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -58129.8
$ws.Range("H141").Value = 108260
$ws.Range("I141").Value = 99890
$ws.Range("J141").Value = 125000
$ws.Range("K141").Value = 99890
$ws.Range("L141").Value = 125000
$ws.Range("M141").Value = -94710
$ws.Range("N141").Value = -135360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2889.9333
$ws.Range("J20").Value = 3654.2856
$ws.Range("L20").Value = 3654.2856
$ws.Range("N20").Value = -4148.2856
$ws.Range("H57").Value = 67154
$ws.Range("J57").Value = 67154
$ws.Range("L57").Value = 67154
$ws.Range("N57").Value = -68594
$ws.Range("H99").Value = 3546.0588
$ws.Range("I99").Value = 3741.1667
$ws.Range("J99").Value = 3077.8
$ws.Range("K99").Value = 3741.1667
$ws.Range("L99").Value = 3077.8
$ws.Range("M99").Value = -2243.1667
$ws.Range("N99").Value = -6073.8
$ws.Range("H133").Value = 75833.336
$ws.Range("J133").Value = 75833.336
$ws.Range("L133").Value = 75833.336
$ws.Range("N133").Value = -85953.336
$ws.Range("H134").Value = 3851.5518
$ws.Range("I134").Value = 2461.7144
$ws.Range("K134").Value = 7385.1432
$ws.Range("M134").Value = -4850.1432
$ws.Range("H135").Value = 49999.2
$ws.Range("J135").Value = 49999.2
$ws.Range("L135").Value = 49999.2
$ws.Range("N135").Value = -60139.2
$ws.Range("H136").Value = 67154
$ws.Range("J136").Value = 67154
$ws.Range("L136").Value = 67154
$ws.Range("N136").Value = -77354
$ws.Range("H137").Value = 62499.5
$ws.Range("J137").Value = 62499.5
$ws.Range("L137").Value = 62499.5
$ws.Range("N137").Value = -72699.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2265.4
$ws.Range("I19").Value = 140.14285
$ws.Range("J19").Value = 4125
$ws.Range("K19").Value = 140.14285
$ws.Range("L19").Value = 4125
$ws.Range("M19").Value = 29.85714999999999
$ws.Range("N19").Value = -4465
$ws.Range("H24").Value = 2265.4
$ws.Range("I24").Value = 140.14285
$ws.Range("J24").Value = 4125
$ws.Range("K24").Value = 140.14285
$ws.Range("L24").Value = 4125
$ws.Range("M24").Value = 29.85714999999999
$ws.Range("N24").Value = -4465
$ws.Range("H31").Value = 3586.4
$ws.Range("I31").Value = 2233.6667
$ws.Range("J31").Value = 8997.333000000001
$ws.Range("K31").Value = 2233.6667
$ws.Range("L31").Value = 8997.333000000001
$ws.Range("M31").Value = -1938.6667
$ws.Range("N31").Value = -9587.333000000001
$ws.Range("H34").Value = 3586.4
$ws.Range("I34").Value = 2233.6667
$ws.Range("J34").Value = 8997.333000000001
$ws.Range("K34").Value = 2233.6667
$ws.Range("L34").Value = 8997.333000000001
$ws.Range("M34").Value = -2031.6667
$ws.Range("N34").Value = -9401.333000000001
$ws.Range("H58").Value = 505359
$ws.Range("I58").Value = 1002228
$ws.Range("J58").Value = 8490
$ws.Range("K58").Value = 1002228
$ws.Range("L58").Value = 8490
$ws.Range("M58").Value = -1002025
$ws.Range("N58").Value = -8896
$ws.Range("H86").Value = 9997.666999999999
$ws.Range("I86").Value = 9994
$ws.Range("K86").Value = 9994
$ws.Range("M86").Value = -8871
$ws.Range("H89").Value = 9997.666999999999
$ws.Range("I89").Value = 9994
$ws.Range("K89").Value = 49970
$ws.Range("M89").Value = -44354
$ws.Range("H99").Value = 561312.3
$ws.Range("I99").Value = 4937
$ws.Range("J99").Value = 839500
$ws.Range("K99").Value = 4937
$ws.Range("L99").Value = 839500
$ws.Range("M99").Value = -3439
$ws.Range("N99").Value = -842496
$ws.Range("H126").Value = 561312.3
$ws.Range("I126").Value = 4937
$ws.Range("J126").Value = 839500
$ws.Range("K126").Value = 14811
$ws.Range("L126").Value = 2518500
$ws.Range("M126").Value = -12341
$ws.Range("N126").Value = -2523440
$ws.Range("H132").Value = 5380.1313
$ws.Range("I132").Value = 4900.3335
$ws.Range("K132").Value = 14701.0005
$ws.Range("M132").Value = -12171.0005
$ws.Range("H136").Value = 505359
$ws.Range("I136").Value = 1002228
$ws.Range("J136").Value = 8490
$ws.Range("K136").Value = 3006684
$ws.Range("L136").Value = 25470
$ws.Range("M136").Value = -3004134
$ws.Range("N136").Value = -30570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 233.375
$ws.Range("J12").Value = 248.86667
$ws.Range("L12").Value = 746.60001
$ws.Range("N12").Value = -1092.60001
$ws.Range("H14").Value = 2075.5454
$ws.Range("I14").Value = 2075.5454
$ws.Range("K14").Value = 6226.6362
$ws.Range("M14").Value = -6053.6362
$ws.Range("H38").Value = 39.833332
$ws.Range("J38").Value = 38
$ws.Range("L38").Value = 114
$ws.Range("N38").Value = -808
$ws.Range("H106").Value = 11949.75
$ws.Range("I106").Value = 5933
$ws.Range("K106").Value = 17799
$ws.Range("M106").Value = -16853
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H122").Value = 125922.875
$ws.Range("I122").Value = 983
$ws.Range("J122").Value = 200886.8
$ws.Range("K122").Value = 8847
$ws.Range("L122").Value = 1807981.2
$ws.Range("M122").Value = -6397
$ws.Range("N122").Value = -1812881.2
$ws.Range("H126").Value = 100
$ws.Range("I126").Value = 100
$ws.Range("K126").Value = 300
$ws.Range("M126").Value = 4640
$ws.Range("H133").Value = 6622.25
$ws.Range("I133").Value = 7996.3335
$ws.Range("K133").Value = 23989.0005
$ws.Range("M133").Value = -18929.0005
$ws.Range("H137").Value = 2403.9167
$ws.Range("I137").Value = 2146
$ws.Range("J137").Value = 2765
$ws.Range("K137").Value = 6438
$ws.Range("L137").Value = 8295
$ws.Range("M137").Value = -1338
$ws.Range("N137").Value = -18495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50005510
$ws.Range("I70").Value = 4666.25
$ws.Range("J70").Value = 62505724
$ws.Range("K70").Value = 4666.25
$ws.Range("L70").Value = 62505724
$ws.Range("M70").Value = -4396.25
$ws.Range("N70").Value = -62506264
$ws.Range("H73").Value = 50005510
$ws.Range("I73").Value = 4666.25
$ws.Range("J73").Value = 62505724
$ws.Range("K73").Value = 4666.25
$ws.Range("L73").Value = 62505724
$ws.Range("M73").Value = -3730.25
$ws.Range("N73").Value = -62507596
$ws.Range("H80").Value = 717092.4
$ws.Range("I80").Value = 558433.75
$ws.Range("J80").Value = 1002677.9
$ws.Range("K80").Value = 558433.75
$ws.Range("L80").Value = 1002677.9
$ws.Range("M80").Value = -557435.75
$ws.Range("N80").Value = -1004673.9
$ws.Range("H83").Value = 717092.4
$ws.Range("I83").Value = 558433.75
$ws.Range("J83").Value = 1002677.9
$ws.Range("K83").Value = 2792168.75
$ws.Range("L83").Value = 5013389.5
$ws.Range("M83").Value = -2787176.75
$ws.Range("N83").Value = -5023373.5
$ws.Range("H105").Value = 45223.668
$ws.Range("J105").Value = 45223.668
$ws.Range("L105").Value = 45223.668
$ws.Range("N105").Value = -52211.668
$ws.Range("H132").Value = 1117046
$ws.Range("I132").Value = 2005160.6
$ws.Range("J132").Value = 6902.75
$ws.Range("K132").Value = 6015481.800000001
$ws.Range("L132").Value = 20708.25
$ws.Range("M132").Value = -6012951.800000001
$ws.Range("N132").Value = -25768.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1527.8572
$ws.Range("J22").Value = 1844
$ws.Range("L22").Value = 1844
$ws.Range("N22").Value = -2434
$ws.Range("H27").Value = 1527.8572
$ws.Range("J27").Value = 1844
$ws.Range("L27").Value = 1844
$ws.Range("N27").Value = -2058
$ws.Range("H40").Value = 3754937.5
$ws.Range("I40").Value = 6003600.5
$ws.Range("J40").Value = 7166.5
$ws.Range("K40").Value = 6003600.5
$ws.Range("L40").Value = 7166.5
$ws.Range("M40").Value = -6003464.5
$ws.Range("N40").Value = -7438.5
$ws.Range("H61").Value = 4595.4165
$ws.Range("I61").Value = 3837.1875
$ws.Range("K61").Value = 3837.1875
$ws.Range("M61").Value = -3635.1875
$ws.Range("H113").Value = 4595.4165
$ws.Range("I113").Value = 3837.1875
$ws.Range("K113").Value = 3837.1875
$ws.Range("M113").Value = -1667.1875
$ws.Range("H122").Value = 1619854.9
$ws.Range("I122").Value = 1432300.6
$ws.Range("J122").Value = 1838668.4
$ws.Range("K122").Value = 4296901.800000001
$ws.Range("L122").Value = 5516005.199999999
$ws.Range("M122").Value = -4294451.800000001
$ws.Range("N122").Value = -5520905.199999999
$ws.Range("H132").Value = 6172.727
$ws.Range("I132").Value = 5500
$ws.Range("J132").Value = 6557.143
$ws.Range("K132").Value = 16500
$ws.Range("L132").Value = 19671.429
$ws.Range("M132").Value = -13970
$ws.Range("N132").Value = -24731.429
$ws.Range("H136").Value = 5025
$ws.Range("I136").Value = 4033.3333
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 12099.9999
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -9549.999899999999
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 826.63635
$ws.Range("I100").Value = 579.6
$ws.Range("J100").Value = 1032.5
$ws.Range("K100").Value = 1159.2
$ws.Range("L100").Value = 2065
$ws.Range("M100").Value = -618.2
$ws.Range("N100").Value = -3147
$ws.Range("H132").Value = 3415.4866
$ws.Range("I132").Value = 2509.2083
$ws.Range("J132").Value = 5088.615
$ws.Range("K132").Value = 7527.624899999999
$ws.Range("L132").Value = 15265.845
$ws.Range("M132").Value = -4997.624899999999
$ws.Range("N132").Value = -20325.845
$ws.Range("H136").Value = 478035.84
$ws.Range("I136").Value = 589632.5
$ws.Range("K136").Value = 1768897.5
$ws.Range("M136").Value = -1766347.5
